# Update the table style used by the "Data Sources from LFX" tables
# (was {48338D13-A9CA-4A3B-89B7-9073ACB4FE37}, now {631E9890-4A1E-43CC-99DB-8398D616D4D7})
# across every slide that contains a table.

$p = $ppt.ActivePresentation

$oldStyleId = "{48338D13-A9CA-4A3B-89B7-9073ACB4FE37}"
$newStyleId = "{631E9890-4A1E-43CC-99DB-8398D616D4D7}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shp = $slide.Shapes.Item($shi)
        if ($shp.HasTable) {
            $tbl = $shp.Table
            if ($tbl.Style -eq $oldStyleId) {
                $tbl.ApplyStyle($newStyleId)
            }
        }
    }
}
